# "Apache POI login test implemented"
# The over-long placeholder test password is replaced with the real
# value used by the login test, and the sheet's active selection is
# moved from I15 to G9 (where the author's cursor ended up after editing).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3 held the placeholder "passwordexcellllllllllllllllllllllllllll" -
# trim it down to the real test password "passwordexcel".
$ws.Range("B3").Value = "passwordexcel"

# Update the active cell/selection shown when the sheet is reopened.
$ws.Range("G9").Select()
